# ---------------------------------------------------------------------------
# paises.xlsx -- "Update countries & provincias Spain"
#
# Refreshes the COVID-19 country snapshot table on sheet "Pais":
#   - bumps the "Datos actualizados" timestamp from 20:04 to 21:04
#   - four countries (Mayotte, Paraguay, Republica del Chad,
#     Islas Virgenes Britanicas) move to a new rank in the list, which
#     pushes the countries that used to occupy those ranks down by one row
#   - refreshed case/recovered/critical/death counts for the countries
#     whose numbers changed between the 20:04 and 21:04 snapshots
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Snapshot timestamp (A1) ----------------------------------------------
$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 9 de Mayo de 2020 a las 21:04'

# --- Mayotte (new row 96) cascades down to Albania at row 102 ---
$ws.Cells.Item(96, 1).Value = 'Mayotte'
$ws.Cells.Item(96, 2).Value = 988
$ws.Cells.Item(96, 3).Value = 134
$ws.Cells.Item(96, 4).Value = 492
$ws.Cells.Item(96, 5).Value = 485
$ws.Cells.Item(96, 6).Value = 7
$ws.Cells.Item(96, 7).Value = 1
$ws.Cells.Item(96, 8).Value = 11

$ws.Cells.Item(97, 1).Value = 'Consejo Danes para los Refugiados'
$ws.Cells.Item(97, 2).Value = 937
$ws.Cells.Item(97, 3).Value = 74
$ws.Cells.Item(97, 4).Value = 130
$ws.Cells.Item(97, 5).Value = 768
$ws.Cells.Item(97, 6).Value = 0
$ws.Cells.Item(97, 7).Value = 3
$ws.Cells.Item(97, 8).Value = 39

$ws.Cells.Item(98, 1).Value = 'Kirguistan'
$ws.Cells.Item(98, 2).Value = 931
$ws.Cells.Item(98, 3).Value = 25
$ws.Cells.Item(98, 4).Value = 658
$ws.Cells.Item(98, 5).Value = 261
$ws.Cells.Item(98, 6).Value = 13
$ws.Cells.Item(98, 8).Value = 12

$ws.Cells.Item(99, 1).Value = 'Letonia'
$ws.Cells.Item(99, 2).Value = 930
$ws.Cells.Item(99, 3).Value = 2
$ws.Cells.Item(99, 4).Value = 464
$ws.Cells.Item(99, 5).Value = 448
$ws.Cells.Item(99, 6).Value = 2
$ws.Cells.Item(99, 7).Value = 0
$ws.Cells.Item(99, 8).Value = 18

$ws.Cells.Item(100, 1).Value = 'Guatemala'
$ws.Cells.Item(100, 2).Value = 900
$ws.Cells.Item(100, 3).Value = 68
$ws.Cells.Item(100, 4).Value = 101
$ws.Cells.Item(100, 5).Value = 775
$ws.Cells.Item(100, 6).Value = 5
$ws.Cells.Item(100, 7).Value = 1
$ws.Cells.Item(100, 8).Value = 24

$ws.Cells.Item(101, 1).Value = 'Republica de Chipre'
$ws.Cells.Item(101, 2).Value = 892
$ws.Cells.Item(101, 3).Value = 1
$ws.Cells.Item(101, 4).Value = 401
$ws.Cells.Item(101, 5).Value = 476
$ws.Cells.Item(101, 6).Value = 10
$ws.Cells.Item(101, 8).Value = 15

$ws.Cells.Item(102, 1).Value = 'Albania'
$ws.Cells.Item(102, 2).Value = 856
$ws.Cells.Item(102, 3).Value = 6
$ws.Cells.Item(102, 4).Value = 627
$ws.Cells.Item(102, 5).Value = 198
$ws.Cells.Item(102, 8).Value = 31

# --- Paraguay (new row 114) cascades down to Guinea-Bisau at row 120 ---
$ws.Cells.Item(114, 1).Value = 'Paraguay'
$ws.Cells.Item(114, 2).Value = 689
$ws.Cells.Item(114, 3).Value = 126
$ws.Cells.Item(114, 4).Value = 155
$ws.Cells.Item(114, 5).Value = 524
$ws.Cells.Item(114, 6).Value = 9
$ws.Cells.Item(114, 7).Value = 0
$ws.Cells.Item(114, 8).Value = 10

$ws.Cells.Item(115, 1).Value = 'Kenia'
$ws.Cells.Item(115, 2).Value = 649
$ws.Cells.Item(115, 3).Value = 28
$ws.Cells.Item(115, 4).Value = 207
$ws.Cells.Item(115, 5).Value = 412
$ws.Cells.Item(115, 6).Value = 1
$ws.Cells.Item(115, 7).Value = 1
$ws.Cells.Item(115, 8).Value = 30

$ws.Cells.Item(116, 1).Value = 'San Marino'
$ws.Cells.Item(116, 2).Value = 637
$ws.Cells.Item(116, 3).Value = 14
$ws.Cells.Item(116, 4).Value = 126
$ws.Cells.Item(116, 5).Value = 470
$ws.Cells.Item(116, 6).Value = 3
$ws.Cells.Item(116, 8).Value = 41

$ws.Cells.Item(117, 1).Value = 'Georgia'
$ws.Cells.Item(117, 2).Value = 626
$ws.Cells.Item(117, 3).Value = 3
$ws.Cells.Item(117, 4).Value = 297
$ws.Cells.Item(117, 5).Value = 319
$ws.Cells.Item(117, 6).Value = 6
$ws.Cells.Item(117, 8).Value = 10

$ws.Cells.Item(118, 1).Value = 'Gabon'
$ws.Cells.Item(118, 2).Value = 620
$ws.Cells.Item(118, 3).Value = 0
$ws.Cells.Item(118, 4).Value = 110
$ws.Cells.Item(118, 5).Value = 502
$ws.Cells.Item(118, 6).Value = 1
$ws.Cells.Item(118, 7).Value = 0
$ws.Cells.Item(118, 8).Value = 8

$ws.Cells.Item(119, 1).Value = 'Tayikistan'
$ws.Cells.Item(119, 2).Value = 612
$ws.Cells.Item(119, 3).Value = 90
$ws.Cells.Item(119, 4).Value = 0
$ws.Cells.Item(119, 5).Value = 592
$ws.Cells.Item(119, 7).Value = 8
$ws.Cells.Item(119, 8).Value = 20

$ws.Cells.Item(120, 1).Value = 'Guinea-Bisau'
$ws.Cells.Item(120, 2).Value = 594
$ws.Cells.Item(120, 4).Value = 25
$ws.Cells.Item(120, 5).Value = 567
$ws.Cells.Item(120, 6).Value = 0
$ws.Cells.Item(120, 8).Value = 2

# --- Republica del Chad (new row 133) cascades down to Ruanda at row 138 ---
$ws.Cells.Item(133, 1).Value = 'Republica del Chad'
$ws.Cells.Item(133, 2).Value = 322
$ws.Cells.Item(133, 3).Value = 62
$ws.Cells.Item(133, 4).Value = 53
$ws.Cells.Item(133, 5).Value = 238
$ws.Cells.Item(133, 7).Value = 3
$ws.Cells.Item(133, 8).Value = 31

$ws.Cells.Item(134, 1).Value = 'Sierra Leona'
$ws.Cells.Item(134, 2).Value = 291
$ws.Cells.Item(134, 3).Value = 34
$ws.Cells.Item(134, 4).Value = 58
$ws.Cells.Item(134, 5).Value = 215
$ws.Cells.Item(134, 6).Value = 0
$ws.Cells.Item(134, 7).Value = 1
$ws.Cells.Item(134, 8).Value = 18

$ws.Cells.Item(135, 1).Value = 'Vietnam'
$ws.Cells.Item(135, 2).Value = 288
$ws.Cells.Item(135, 3).Value = 0
$ws.Cells.Item(135, 4).Value = 241
$ws.Cells.Item(135, 5).Value = 47
$ws.Cells.Item(135, 6).Value = 8
$ws.Cells.Item(135, 8).Value = 0

$ws.Cells.Item(136, 1).Value = 'Benin'
$ws.Cells.Item(136, 2).Value = 284
$ws.Cells.Item(136, 3).Value = 42
$ws.Cells.Item(136, 4).Value = 62
$ws.Cells.Item(136, 5).Value = 220
$ws.Cells.Item(136, 8).Value = 2

$ws.Cells.Item(137, 1).Value = 'Congo'
$ws.Cells.Item(137, 2).Value = 274
$ws.Cells.Item(137, 4).Value = 33
$ws.Cells.Item(137, 5).Value = 231
$ws.Cells.Item(137, 8).Value = 10

$ws.Cells.Item(138, 1).Value = 'Ruanda'
$ws.Cells.Item(138, 2).Value = 273
$ws.Cells.Item(138, 4).Value = 136
$ws.Cells.Item(138, 5).Value = 137
$ws.Cells.Item(138, 8).Value = 0

# --- Islas Virgenes Britanicas / Butan swap at rows 212-213 ---
$ws.Cells.Item(212, 1).Value = 'Islas Virgenes Britanicas'
$ws.Cells.Item(212, 4).Value = 4
$ws.Cells.Item(212, 8).Value = 1

$ws.Cells.Item(213, 1).Value = 'Butan'
$ws.Cells.Item(213, 4).Value = 5
$ws.Cells.Item(213, 8).Value = 0

# --- Updated counts for countries whose rank did not change ---------------
$ws.Cells.Item(4, 2).Value = 1335668
$ws.Cells.Item(4, 3).Value = 13883
$ws.Cells.Item(4, 5).Value = 1031567
$ws.Cells.Item(4, 7).Value = 853
$ws.Cells.Item(4, 8).Value = 79468

$ws.Cells.Item(9, 4).Value = 56038
$ws.Cells.Item(9, 5).Value = 93731
$ws.Cells.Item(9, 6).Value = 2812
$ws.Cells.Item(9, 7).Value = 80
$ws.Cells.Item(9, 8).Value = 26310

$ws.Cells.Item(15, 2).Value = 67619
$ws.Cells.Item(15, 3).Value = 1185
$ws.Cells.Item(15, 4).Value = 30980
$ws.Cells.Item(15, 5).Value = 31949
$ws.Cells.Item(15, 7).Value = 121
$ws.Cells.Item(15, 8).Value = 4690

$ws.Cells.Item(16, 2).Value = 62808
$ws.Cells.Item(16, 3).Value = 3113
$ws.Cells.Item(16, 4).Value = 19301
$ws.Cells.Item(16, 5).Value = 41406
$ws.Cells.Item(16, 7).Value = 116
$ws.Cells.Item(16, 8).Value = 2101

$ws.Cells.Item(23, 2).Value = 29071
$ws.Cells.Item(23, 3).Value = 253
$ws.Cells.Item(23, 5).Value = 23921
$ws.Cells.Item(23, 6).Value = 181
$ws.Cells.Item(23, 7).Value = 13
$ws.Cells.Item(23, 8).Value = 1717

$ws.Cells.Item(50, 2).Value = 8084
$ws.Cells.Item(50, 3).Value = 14
$ws.Cells.Item(50, 5).Value = 7833
$ws.Cells.Item(50, 6).Value = 24
$ws.Cells.Item(50, 7).Value = 1
$ws.Cells.Item(50, 8).Value = 219

$ws.Cells.Item(60, 2).Value = 4867
$ws.Cells.Item(60, 3).Value = 139
$ws.Cells.Item(60, 5).Value = 2781
$ws.Cells.Item(60, 7).Value = 11
$ws.Cells.Item(60, 8).Value = 161

$ws.Cells.Item(110, 2).Value = 748
$ws.Cells.Item(110, 3).Value = 4
$ws.Cells.Item(110, 4).Value = 569
$ws.Cells.Item(110, 5).Value = 131

$ws.Cells.Item(127, 2).Value = 431
$ws.Cells.Item(127, 3).Value = 5
$ws.Cells.Item(127, 5).Value = 77
$ws.Cells.Item(127, 6).Value = 5

$ws.Cells.Item(140, 2).Value = 236
$ws.Cells.Item(140, 3).Value = 6
$ws.Cells.Item(140, 5).Value = 190

$ws.Cells.Item(176, 4).Value = 13
$ws.Cells.Item(176, 5).Value = 28

$ws.Cells.Item(186, 4).Value = 12
$ws.Cells.Item(186, 5).Value = 10
